$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.2120343839541547
$ws.Range("C2").Value = 0.5071633237822349
$ws.Range("J2").Value = 0.01432664756446991
$ws.Range("P2").Value = 0.1776504297994269
$ws.Range("S2").Value = 0.08882521489971347

$ws.Range("B3").Value = 0.0111731843575419
$ws.Range("C3").Value = 0.00558659217877095
$ws.Range("J3").Value = 0.0446927374301676
$ws.Range("P3").Value = 0.6815642458100558
$ws.Range("S3").Value = 0.2569832402234637

$ws.Range("J4").Value = 0.03508771929824561
$ws.Range("P4").Value = 0.7543859649122807
$ws.Range("S4").Value = 0.2105263157894737

$ws.Range("B6").Value = 0.04564315352697095
$ws.Range("D6").Value = 0.02489626556016597
$ws.Range("F6").Value = 0.1161825726141079
$ws.Range("J6").Value = 0.2365145228215768
$ws.Range("O6").Value = 0.02904564315352697
$ws.Range("Q6").Value = 0.1410788381742739
$ws.Range("R6").Value = 0.07468879668049792
$ws.Range("S6").Value = 0.3319502074688797

$ws.Range("B7").Value = 0.109704641350211
$ws.Range("D7").Value = 0.02109704641350211
$ws.Range("F7").Value = 0.0379746835443038
$ws.Range("J7").Value = 0.189873417721519
$ws.Range("O7").Value = 0.01265822784810127
$ws.Range("Q7").Value = 0.160337552742616
$ws.Range("R7").Value = 0.08860759493670886
$ws.Range("S7").Value = 0.379746835443038

$ws.Range("B8").Value = 0.09380234505862646
$ws.Range("D8").Value = 0.01842546063651591
$ws.Range("F8").Value = 0.05862646566164154
$ws.Range("J8").Value = 0.1289782244556114
$ws.Range("O8").Value = 0.01507537688442211
$ws.Range("Q8").Value = 0.1775544388609715
$ws.Range("R8").Value = 0.1122278056951424
$ws.Range("S8").Value = 0.3953098827470687

$ws.Range("B9").Value = 0.096
$ws.Range("D9").Value = 0.012
$ws.Range("F9").Value = 0.064
$ws.Range("J9").Value = 0.104
$ws.Range("O9").Value = 0.004
$ws.Range("Q9").Value = 0.2
$ws.Range("R9").Value = 0.076
$ws.Range("S9").Value = 0.444

$ws.Range("B10").Value = 0.1097560975609756
$ws.Range("D10").Value = 0.02510760401721664
$ws.Range("E10").Value = 0.0007173601147776184
$ws.Range("F10").Value = 0.054519368723099
$ws.Range("J10").Value = 0.1248206599713056
$ws.Range("O10").Value = 0.01506456241032999
$ws.Range("Q10").Value = 0.2022955523672884
$ws.Range("R10").Value = 0.07962697274031563
$ws.Range("S10").Value = 0.3880918220946915

$ws.Range("G11").Value = 0.1420118343195266
$ws.Range("J11").Value = 0.09763313609467456
$ws.Range("K11").Value = 0.1834319526627219
$ws.Range("L11").Value = 0.5769230769230769

$ws.Range("G12").Value = 0.75
$ws.Range("J12").Value = 0.22
$ws.Range("L12").Value = 0.02
$ws.Range("S12").Value = 0.01

$ws.Range("G13").Value = 0.75
$ws.Range("J13").Value = 0.2333333333333333
$ws.Range("S13").Value = 0.01666666666666667

$ws.Range("F15").Value = 0.02109704641350211
$ws.Range("H15").Value = 0.1687763713080169
$ws.Range("I15").Value = 0.0759493670886076
$ws.Range("J15").Value = 0.3375527426160337
$ws.Range("K15").Value = 0.0379746835443038
$ws.Range("O15").Value = 0.08016877637130802
$ws.Range("S15").Value = 0.2784810126582278

$ws.Range("F16").Value = 0.02262443438914027
$ws.Range("H16").Value = 0.1900452488687783
$ws.Range("I16").Value = 0.08144796380090498
$ws.Range("J16").Value = 0.4072398190045249
$ws.Range("K16").Value = 0.1357466063348416
$ws.Range("M16").Value = 0.02262443438914027
$ws.Range("N16").Value = 0.004524886877828055
$ws.Range("O16").Value = 0.04977375565610859
$ws.Range("S16").Value = 0.08597285067873303

$ws.Range("F17").Value = 0.01568627450980392
$ws.Range("H17").Value = 0.1921568627450981
$ws.Range("I17").Value = 0.1098039215686274
$ws.Range("J17").Value = 0.396078431372549
$ws.Range("K17").Value = 0.09607843137254903
$ws.Range("M17").Value = 0.02941176470588235
$ws.Range("O17").Value = 0.05490196078431372
$ws.Range("S17").Value = 0.1058823529411765

$ws.Range("F18").Value = 0.02155172413793104
$ws.Range("H18").Value = 0.2241379310344828
$ws.Range("I18").Value = 0.06896551724137931
$ws.Range("J18").Value = 0.3879310344827586
$ws.Range("K18").Value = 0.1293103448275862
$ws.Range("M18").Value = 0.02155172413793104
$ws.Range("O18").Value = 0.05172413793103448
$ws.Range("S18").Value = 0.09482758620689655

$ws.Range("F19").Value = 0.01783264746227709
$ws.Range("H19").Value = 0.2517146776406036
$ws.Range("I19").Value = 0.1008230452674897
$ws.Range("J19").Value = 0.3209876543209876
$ws.Range("K19").Value = 0.1042524005486968
$ws.Range("M19").Value = 0.02400548696844993
$ws.Range("O19").Value = 0.06721536351165981
$ws.Range("S19").Value = 0.1131687242798354
